# Scheduled-runner style update of market/price-derived figures across the
# Titan_Profits leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Only the price/profit columns (H:N) of specific rows are refreshed; all
# other data (leve names, items, levels, etc.) is left untouched.

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H135").Value = 13080.333
$ws.Range("I135").Value = 13080.333
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 117722.997
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -115187.997

$ws.Range("H137").Value = 28572936
$ws.Range("I137").Value = 37038410
$ws.Range("J137").Value = 1963.25
$ws.Range("K137").Value = 111115230
$ws.Range("L137").Value = 5889.75
$ws.Range("M137").Value = -111112680
$ws.Range("N137").Value = -10989.75

$ws.Range("H138").Value = 5489397
$ws.Range("I138").Value = 1241513.2
$ws.Range("J138").Value = 7044960
$ws.Range("K138").Value = 3724539.6
$ws.Range("L138").Value = 21134880
$ws.Range("M138").Value = -3719399.6
$ws.Range("N138").Value = -21145160

# --- ARM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H3").Value = 1750
$ws.Range("I3").Value = 1333.3334
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 1333.3334
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -1218.3334
$ws.Range("N3").Value = -3230

$ws.Range("H32").Value = 21874.889
$ws.Range("I32").Value = 2482.0293
$ws.Range("J32").Value = 54842.75
$ws.Range("K32").Value = 2482.0293
$ws.Range("L32").Value = 54842.75
$ws.Range("M32").Value = -2195.0293
$ws.Range("N32").Value = -55416.75

$ws.Range("H61").Value = 2721.5186
$ws.Range("I61").Value = 2111.25
$ws.Range("J61").Value = 4465.143
$ws.Range("K61").Value = 2111.25
$ws.Range("L61").Value = 4465.143
$ws.Range("M61").Value = -1899.25
$ws.Range("N61").Value = -4889.143

$ws.Range("H74").Value = 7581.15
$ws.Range("I74").Value = 1821.6
$ws.Range("J74").Value = 24859.8
$ws.Range("K74").Value = 1821.6
$ws.Range("L74").Value = 24859.8
$ws.Range("M74").Value = -947.5999999999999
$ws.Range("N74").Value = -26607.8

$ws.Range("H77").Value = 7581.15
$ws.Range("I77").Value = 1821.6
$ws.Range("J77").Value = 24859.8
$ws.Range("K77").Value = 9108
$ws.Range("L77").Value = 124299
$ws.Range("M77").Value = -4740
$ws.Range("N77").Value = -133035

$ws.Range("H132").Value = 1714.5
$ws.Range("I132").Value = 1392.2709
$ws.Range("J132").Value = 4292.3335
$ws.Range("K132").Value = 4176.8127
$ws.Range("L132").Value = 12877.0005
$ws.Range("M132").Value = -1646.8127
$ws.Range("N132").Value = -17937.0005

$ws.Range("H136").Value = 2721.5186
$ws.Range("I136").Value = 2111.25
$ws.Range("J136").Value = 4465.143
$ws.Range("K136").Value = 6333.75
$ws.Range("L136").Value = 13395.429
$ws.Range("M136").Value = -3783.75
$ws.Range("N136").Value = -18495.429

# --- BSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H94").Value = 1842.4166
$ws.Range("I94").Value = 1842.4166
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1842.4166
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1391.4166

$ws.Range("H107").Value = 605.2
$ws.Range("I107").Value = 588.2941
$ws.Range("J107").Value = 641.125
$ws.Range("K107").Value = 588.2941
$ws.Range("L107").Value = 641.125
$ws.Range("M107").Value = 1331.7059
$ws.Range("N107").Value = -4481.125

$ws.Range("H134").Value = 5002.1055
$ws.Range("I134").Value = 3411
$ws.Range("J134").Value = 6770
$ws.Range("K134").Value = 10233
$ws.Range("L134").Value = 20310
$ws.Range("M134").Value = -7698
$ws.Range("N134").Value = -25380

# --- CRP -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H21").Value = 6453.75
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 6453.75
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 6453.75
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -6923.75

$ws.Range("H99").Value = 15626081
$ws.Range("I99").Value = 15626081
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 15626081
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -15624583

$ws.Range("H126").Value = 15626081
$ws.Range("I126").Value = 15626081
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 46878243
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -46875773

$ws.Range("H132").Value = 3064.9546
$ws.Range("I132").Value = 2372.6
$ws.Range("J132").Value = 4548.5713
$ws.Range("K132").Value = 7117.799999999999
$ws.Range("L132").Value = 13645.7139
$ws.Range("M132").Value = -4587.799999999999
$ws.Range("N132").Value = -18705.7139

# --- CUL ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H12").Value = 36
$ws.Range("I12").Value = 13.714286
$ws.Range("J12").Value = 75
$ws.Range("K12").Value = 41.142858
$ws.Range("L12").Value = 225
$ws.Range("M12").Value = 131.857142

$ws.Range("H109").Value = 2750
$ws.Range("I109").Value = 2000
$ws.Range("J109").Value = 3500
$ws.Range("K109").Value = 6000
$ws.Range("L109").Value = 10500
$ws.Range("M109").Value = -4960
$ws.Range("N109").Value = -12580

$ws.Range("H131").Value = 1503.3103
$ws.Range("I131").Value = 449.5
$ws.Range("J131").Value = 1722.8541
$ws.Range("K131").Value = 1348.5
$ws.Range("L131").Value = 5168.5623
$ws.Range("M131").Value = 3691.5
$ws.Range("N131").Value = -15248.5623

# --- GSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H57").Value = 22500
$ws.Range("I57").Value = 22500
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 22500
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -21680

$ws.Range("H97").Value = 839.3333
$ws.Range("I97").Value = 770.7143
$ws.Range("J97").Value = 1800
$ws.Range("K97").Value = 770.7143
$ws.Range("L97").Value = 1800
$ws.Range("M97").Value = -274.7143

$ws.Range("H107").Value = 450.125
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 480.2
$ws.Range("K107").Value = 400
$ws.Range("L107").Value = 480.2
$ws.Range("M107").Value = 1520
$ws.Range("N107").Value = -4320.2

$ws.Range("H132").Value = 2455.878
$ws.Range("I132").Value = 1963.5294
$ws.Range("J132").Value = 4847.2856
$ws.Range("K132").Value = 5890.5882
$ws.Range("L132").Value = 14541.8568
$ws.Range("M132").Value = -3360.5882
$ws.Range("N132").Value = -19601.8568

# --- LTW ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H93").Value = 506.85715
$ws.Range("I93").Value = 539.38464
$ws.Range("J93").Value = 84
$ws.Range("K93").Value = 539.38464
$ws.Range("L93").Value = 84
$ws.Range("M93").Value = 708.61536
$ws.Range("N93").Value = -2580

$ws.Range("H94").Value = 21666.666
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 21666.666
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 21666.666
$ws.Range("N94").Value = -23018.666

# --- WVR ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H2").Value = 15498.733
$ws.Range("I2").Value = 50000
$ws.Range("J2").Value = 13034.357
$ws.Range("K2").Value = 50000
$ws.Range("L2").Value = 13034.357
$ws.Range("M2").Value = -49888
$ws.Range("N2").Value = -13258.357

$ws.Range("H4").Value = 21937.875
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 21937.875
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 21937.875
$ws.Range("N4").Value = -22163.875

$ws.Range("H136").Value = 18575644
$ws.Range("I136").Value = 25719018
$ws.Range("J136").Value = 2875
$ws.Range("K136").Value = 77157054
$ws.Range("L136").Value = 8625
$ws.Range("M136").Value = -77154504
$ws.Range("N136").Value = -13725
